$d = $word.ActiveDocument

# 1. Insert a new paragraph BEFORE "Migliorare posizione del cestino..." paragraph.
$r = $d.Content.Find.Execute("Migliorare posizione del cestino per la rimozione su certe notifiche;", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Migliorare posizione del cestino per la rimozione su certe notifiche;*") {
        $target = $p
        break
    }
}

$insertRange = $target.Range
$insertRange.Collapse(1)  # wdCollapseStart
$insertRange.InsertBefore("Aggiungere spazio alla fine dei grafici; PARIX`r")

